$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "4+14="
$t.Rows.Item(1).Cells.Item(2).Range.Text = "9+37="
$t.Rows.Item(1).Cells.Item(3).Range.Text = "82-9="
$t.Rows.Item(1).Cells.Item(4).Range.Text = "89-2="
$t.Rows.Item(1).Cells.Item(5).Range.Text = "88+9="
$t.Rows.Item(2).Cells.Item(1).Range.Text = "95-5="
$t.Rows.Item(2).Cells.Item(2).Range.Text = "14-13="
$t.Rows.Item(2).Cells.Item(3).Range.Text = "37-32="
$t.Rows.Item(2).Cells.Item(4).Range.Text = "15-7="
$t.Rows.Item(2).Cells.Item(5).Range.Text = "79+14="
$t.Rows.Item(3).Cells.Item(1).Range.Text = "84-59="
$t.Rows.Item(3).Cells.Item(2).Range.Text = "10+13="
$t.Rows.Item(3).Cells.Item(3).Range.Text = "13+50="
$t.Rows.Item(3).Cells.Item(4).Range.Text = "36+38="
$t.Rows.Item(3).Cells.Item(5).Range.Text = "20+51="
$t.Rows.Item(4).Cells.Item(1).Range.Text = "50+35="
$t.Rows.Item(4).Cells.Item(2).Range.Text = "56-10="
$t.Rows.Item(4).Cells.Item(3).Range.Text = "57+0="
$t.Rows.Item(4).Cells.Item(4).Range.Text = "69-0="
$t.Rows.Item(4).Cells.Item(5).Range.Text = "8+50="
$t.Rows.Item(5).Cells.Item(1).Range.Text = "11+54="
$t.Rows.Item(5).Cells.Item(2).Range.Text = "3+13="
$t.Rows.Item(5).Cells.Item(3).Range.Text = "94-87="
$t.Rows.Item(5).Cells.Item(4).Range.Text = "19-19="
$t.Rows.Item(5).Cells.Item(5).Range.Text = "46-10="
$t.Rows.Item(6).Cells.Item(1).Range.Text = "2+60="
$t.Rows.Item(6).Cells.Item(2).Range.Text = "55+38="
$t.Rows.Item(6).Cells.Item(3).Range.Text = "32-16="
$t.Rows.Item(6).Cells.Item(4).Range.Text = "25-18="
$t.Rows.Item(6).Cells.Item(5).Range.Text = "57-55="
$t.Rows.Item(7).Cells.Item(1).Range.Text = "3+19="
$t.Rows.Item(7).Cells.Item(2).Range.Text = "21+33="
$t.Rows.Item(7).Cells.Item(3).Range.Text = "40+18="
$t.Rows.Item(7).Cells.Item(4).Range.Text = "74-57="
$t.Rows.Item(7).Cells.Item(5).Range.Text = "46-36="
$t.Rows.Item(8).Cells.Item(1).Range.Text = "42-30="
$t.Rows.Item(8).Cells.Item(2).Range.Text = "78-41="
$t.Rows.Item(8).Cells.Item(3).Range.Text = "79-23="
$t.Rows.Item(8).Cells.Item(4).Range.Text = "10+59="
$t.Rows.Item(8).Cells.Item(5).Range.Text = "74+8="
$t.Rows.Item(9).Cells.Item(1).Range.Text = "75-32="
$t.Rows.Item(9).Cells.Item(2).Range.Text = "80+9="
$t.Rows.Item(9).Cells.Item(3).Range.Text = "51+40="
$t.Rows.Item(9).Cells.Item(4).Range.Text = "11+69="
$t.Rows.Item(9).Cells.Item(5).Range.Text = "2+45="
$t.Rows.Item(10).Cells.Item(1).Range.Text = "20-4="
$t.Rows.Item(10).Cells.Item(2).Range.Text = "36-30="
$t.Rows.Item(10).Cells.Item(3).Range.Text = "45+42="
$t.Rows.Item(10).Cells.Item(4).Range.Text = "54-49="
$t.Rows.Item(10).Cells.Item(5).Range.Text = "78-20="
$t.Rows.Item(11).Cells.Item(1).Range.Text = "32-18="
$t.Rows.Item(11).Cells.Item(2).Range.Text = "64-25="
$t.Rows.Item(11).Cells.Item(3).Range.Text = "24+66="
$t.Rows.Item(11).Cells.Item(4).Range.Text = "1+34="
$t.Rows.Item(11).Cells.Item(5).Range.Text = "93-77="
$t.Rows.Item(12).Cells.Item(1).Range.Text = "76-27="
$t.Rows.Item(12).Cells.Item(2).Range.Text = "8+83="
$t.Rows.Item(12).Cells.Item(3).Range.Text = "20+22="
$t.Rows.Item(12).Cells.Item(4).Range.Text = "95-39="
$t.Rows.Item(12).Cells.Item(5).Range.Text = "48-41="
$t.Rows.Item(13).Cells.Item(1).Range.Text = "49+46="
$t.Rows.Item(13).Cells.Item(2).Range.Text = "37-22="
$t.Rows.Item(13).Cells.Item(3).Range.Text = "79-53="
$t.Rows.Item(13).Cells.Item(4).Range.Text = "61-27="
$t.Rows.Item(13).Cells.Item(5).Range.Text = "85-63="
$t.Rows.Item(14).Cells.Item(1).Range.Text = "86-27="
$t.Rows.Item(14).Cells.Item(2).Range.Text = "91-47="
$t.Rows.Item(14).Cells.Item(3).Range.Text = "6+16="
$t.Rows.Item(14).Cells.Item(4).Range.Text = "72-17="
$t.Rows.Item(14).Cells.Item(5).Range.Text = "95-52="
$t.Rows.Item(15).Cells.Item(1).Range.Text = "37-12="
$t.Rows.Item(15).Cells.Item(2).Range.Text = "72-8="
$t.Rows.Item(15).Cells.Item(3).Range.Text = "76+6="
$t.Rows.Item(15).Cells.Item(4).Range.Text = "81-1="
$t.Rows.Item(15).Cells.Item(5).Range.Text = "50-27="
$t.Rows.Item(16).Cells.Item(1).Range.Text = "60+20="
$t.Rows.Item(16).Cells.Item(2).Range.Text = "51+33="
$t.Rows.Item(16).Cells.Item(3).Range.Text = "6+86="
$t.Rows.Item(16).Cells.Item(4).Range.Text = "29+66="
$t.Rows.Item(16).Cells.Item(5).Range.Text = "16+50="
$t.Rows.Item(17).Cells.Item(1).Range.Text = "96-41="
$t.Rows.Item(17).Cells.Item(2).Range.Text = "86-56="
$t.Rows.Item(17).Cells.Item(3).Range.Text = "99-34="
$t.Rows.Item(17).Cells.Item(4).Range.Text = "83-62="
$t.Rows.Item(17).Cells.Item(5).Range.Text = "57-13="
$t.Rows.Item(18).Cells.Item(1).Range.Text = "84-4="
$t.Rows.Item(18).Cells.Item(2).Range.Text = "26+17="
$t.Rows.Item(18).Cells.Item(3).Range.Text = "25+11="
$t.Rows.Item(18).Cells.Item(4).Range.Text = "81-29="
$t.Rows.Item(18).Cells.Item(5).Range.Text = "84-47="
$t.Rows.Item(19).Cells.Item(1).Range.Text = "19+40="
$t.Rows.Item(19).Cells.Item(2).Range.Text = "73+4="
$t.Rows.Item(19).Cells.Item(3).Range.Text = "28+60="
$t.Rows.Item(19).Cells.Item(4).Range.Text = "83+8="
$t.Rows.Item(19).Cells.Item(5).Range.Text = "76+11="
$t.Rows.Item(20).Cells.Item(1).Range.Text = "81-33="
$t.Rows.Item(20).Cells.Item(2).Range.Text = "40-33="
$t.Rows.Item(20).Cells.Item(3).Range.Text = "6+5="
$t.Rows.Item(20).Cells.Item(4).Range.Text = "10+3="
$t.Rows.Item(20).Cells.Item(5).Range.Text = "76-7="
